$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column for all existing data rows
# (rows 2-27) from 2023-12-07 (45267) to 2023-12-08 (45268).
$ws.Range("C2:C27").Value = 45268

# Row 27 gains an explicit custom row height marker in the saved file.
$ws.Rows.Item(27).RowHeight = 15

# Append the new data row (row 28) describing case "A 62339-2023".
$ws.Range("A28").Value = "A 62339-2023"

$ws.Range("B28").Value = 45267
$ws.Range("B28").NumberFormat = "YYYY-MM-DD"

$ws.Range("C28").Value = 45268
$ws.Range("C28").NumberFormat = "YYYY-MM-DD"

$ws.Range("D28").Value = "OKÄNT"
$ws.Range("E28").Value = "OKÄNT"

$ws.Range("G28").Value = 2.9
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = 0
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = 0
$ws.Range("Q28").Value = 0

$ws.Range("R28").WrapText = $true
